$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 9109.333000000001  # ALC!H40 was 8682.474
$ws.Cells.Item(40, 9).Value = 11689  # ALC!I40 was 10798.167
$ws.Cells.Item(40, 11).Value = 11689  # ALC!K40 was 10798.167
$ws.Cells.Item(40, 13).Value = -11514  # ALC!M40 was -10623.167

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(61, 8).Value = 587  # ALC!H61 was 581
$ws.Cells.Item(61, 9).Value = 587  # ALC!I61 was 571.25
$ws.Cells.Item(61, 10).Value = 0  # ALC!J61 was 620
$ws.Cells.Item(61, 11).Value = 1761  # ALC!K61 was 1713.75
$ws.Cells.Item(61, 12).Value = 0  # ALC!L61 was 1860
$ws.Cells.Item(61, 13).Value = -1589  # ALC!M61 was -1541.75
$ws.Cells.Item(61, 14).ClearContents()  # ALC!N61 was -2204

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 4333.6665  # ALC!H64 was 4750.6665
$ws.Cells.Item(64, 9).Value = 3750.5  # ALC!I64 was 3501
$ws.Cells.Item(64, 10).Value = 5500  # ALC!J64 was 7250
$ws.Cells.Item(64, 11).Value = 3750.5  # ALC!K64 was 3501
$ws.Cells.Item(64, 12).Value = 5500  # ALC!L64 was 7250
$ws.Cells.Item(64, 13).Value = -3502.5  # ALC!M64 was -3253
$ws.Cells.Item(64, 14).Value = -5996  # ALC!N64 was -7746

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(67, 8).Value = 4333.6665  # ALC!H67 was 4750.6665
$ws.Cells.Item(67, 9).Value = 3750.5  # ALC!I67 was 3501
$ws.Cells.Item(67, 10).Value = 5500  # ALC!J67 was 7250
$ws.Cells.Item(67, 11).Value = 3750.5  # ALC!K67 was 3501
$ws.Cells.Item(67, 12).Value = 5500  # ALC!L67 was 7250
$ws.Cells.Item(67, 13).Value = -2892.5  # ALC!M67 was -2643
$ws.Cells.Item(67, 14).Value = -7216  # ALC!N67 was -8966

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 9041.5  # ALC!H74 was 8744.941000000001
$ws.Cells.Item(74, 9).Value = 4189  # ALC!I74 was 4151.2
$ws.Cells.Item(74, 11).Value = 4189  # ALC!K74 was 4151.2
$ws.Cells.Item(74, 13).Value = -3253  # ALC!M74 was -3215.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(77, 8).Value = 9041.5  # ALC!H77 was 8744.941000000001
$ws.Cells.Item(77, 9).Value = 4189  # ALC!I77 was 4151.2
$ws.Cells.Item(77, 11).Value = 20945  # ALC!K77 was 20756
$ws.Cells.Item(77, 13).Value = -16265  # ALC!M77 was -16076

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(87, 8).Value = 59991.5  # ALC!H87 was 59992
$ws.Cells.Item(87, 10).Value = 59991.5  # ALC!J87 was 59992
$ws.Cells.Item(87, 12).Value = 59991.5  # ALC!L87 was 59992
$ws.Cells.Item(87, 14).Value = -62487.5  # ALC!N87 was -62488

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(90, 8).Value = 59991.5  # ALC!H90 was 59992
$ws.Cells.Item(90, 10).Value = 59991.5  # ALC!J90 was 59992
$ws.Cells.Item(90, 12).Value = 179974.5  # ALC!L90 was 179976
$ws.Cells.Item(90, 14).Value = -192454.5  # ALC!N90 was -192456

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(24, 8).Value = 54797  # ARM!H24 was 54576
$ws.Cells.Item(24, 10).Value = 54797  # ARM!J24 was 54576
$ws.Cells.Item(24, 12).Value = 54797  # ARM!L24 was 54576
$ws.Cells.Item(24, 14).Value = -55545  # ARM!N24 was -55324

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1813708.9  # ARM!H32 was 1867849.4
$ws.Cells.Item(32, 9).Value = 1840307.5  # ARM!I32 was 1896074.4
$ws.Cells.Item(32, 11).Value = 1840307.5  # ARM!K32 was 1896074.4
$ws.Cells.Item(32, 13).Value = -1840020.5  # ARM!M32 was -1895787.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 993.5  # ARM!H63 was 1247
$ws.Cells.Item(63, 9).Value = 993.5  # ARM!I63 was 994
$ws.Cells.Item(63, 10).Value = 0  # ARM!J63 was 1500
$ws.Cells.Item(63, 11).Value = 993.5  # ARM!K63 was 994
$ws.Cells.Item(63, 12).Value = 0  # ARM!L63 was 1500
$ws.Cells.Item(63, 13).Value = -307.5  # ARM!M63 was -308
$ws.Cells.Item(63, 14).ClearContents()  # ARM!N63 was -2872

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value = 993.5  # ARM!H66 was 1247
$ws.Cells.Item(66, 9).Value = 993.5  # ARM!I66 was 994
$ws.Cells.Item(66, 10).Value = 0  # ARM!J66 was 1500
$ws.Cells.Item(66, 11).Value = 4967.5  # ARM!K66 was 4970
$ws.Cells.Item(66, 12).Value = 0  # ARM!L66 was 7500
$ws.Cells.Item(66, 13).Value = -1535.5  # ARM!M66 was -1538
$ws.Cells.Item(66, 14).ClearContents()  # ARM!N66 was -14364

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(80, 8).Value = 50422  # ARM!H80 was 46357.6
$ws.Cells.Item(80, 9).Value = 30000  # ARM!I80 was 30050
$ws.Cells.Item(80, 11).Value = 30000  # ARM!K80 was 30050
$ws.Cells.Item(80, 13).Value = -29002  # ARM!M80 was -29052

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(83, 8).Value = 50422  # ARM!H83 was 46357.6
$ws.Cells.Item(83, 9).Value = 30000  # ARM!I83 was 30050
$ws.Cells.Item(83, 11).Value = 90000  # ARM!K83 was 90150
$ws.Cells.Item(83, 13).Value = -85008  # ARM!M83 was -85158

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(100, 8).Value = 54797  # ARM!H100 was 54576
$ws.Cells.Item(100, 10).Value = 54797  # ARM!J100 was 54576
$ws.Cells.Item(100, 12).Value = 54797  # ARM!L100 was 54576
$ws.Cells.Item(100, 14).Value = -56961  # ARM!N100 was -56740

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 22224094  # ARM!H110 was 16668479
$ws.Cells.Item(110, 9).Value = 1927.3077  # ARM!I110 was 1846.4445
$ws.Cells.Item(110, 11).Value = 1927.3077  # ARM!K110 was 1846.4445
$ws.Cells.Item(110, 13).Value = 117.6922999999999  # ARM!M110 was 198.5554999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 8154.486  # ARM!H132 was 8154.514
$ws.Cells.Item(132, 9).Value = 8355.375  # ARM!I132 was 8355.4375
$ws.Cells.Item(132, 11).Value = 25066.125  # ARM!K132 was 25066.3125
$ws.Cells.Item(132, 13).Value = -22536.125  # ARM!M132 was -22536.3125

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 35714612  # BSM!H80 was 33333652
$ws.Cells.Item(80, 10).Value = 276.5  # BSM!J80 was 272.9091
$ws.Cells.Item(80, 12).Value = 276.5  # BSM!L80 was 272.9091
$ws.Cells.Item(80, 14).Value = -2272.5  # BSM!N80 was -2268.9091

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82, 8).Value = 4367.4  # BSM!H82 was 15778.6
$ws.Cells.Item(82, 9).Value = 4367.4  # BSM!I82 was 4723.25
$ws.Cells.Item(82, 10).Value = 0  # BSM!J82 was 60000
$ws.Cells.Item(82, 11).Value = 4367.4  # BSM!K82 was 4723.25
$ws.Cells.Item(82, 12).Value = 0  # BSM!L82 was 60000
$ws.Cells.Item(82, 13).Value = -3984.4  # BSM!M82 was -4340.25
$ws.Cells.Item(82, 14).ClearContents()  # BSM!N82 was -60766

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(83, 8).Value = 35714612  # BSM!H83 was 33333652
$ws.Cells.Item(83, 10).Value = 276.5  # BSM!J83 was 272.9091
$ws.Cells.Item(83, 12).Value = 1382.5  # BSM!L83 was 1364.5455
$ws.Cells.Item(83, 14).Value = -11366.5  # BSM!N83 was -11348.5455

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(85, 8).Value = 4367.4  # BSM!H85 was 15778.6
$ws.Cells.Item(85, 9).Value = 4367.4  # BSM!I85 was 4723.25
$ws.Cells.Item(85, 10).Value = 0  # BSM!J85 was 60000
$ws.Cells.Item(85, 11).Value = 4367.4  # BSM!K85 was 4723.25
$ws.Cells.Item(85, 12).Value = 0  # BSM!L85 was 60000
$ws.Cells.Item(85, 13).Value = -3041.4  # BSM!M85 was -3397.25
$ws.Cells.Item(85, 14).ClearContents()  # BSM!N85 was -62652

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 15153514  # BSM!H99 was 18184188
$ws.Cells.Item(99, 9).Value = 1874.75  # BSM!I99 was 2450
$ws.Cells.Item(99, 11).Value = 1874.75  # BSM!K99 was 2450
$ws.Cells.Item(99, 13).Value = -376.75  # BSM!M99 was -952

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 3502.348  # BSM!H105 was 3511.4348
$ws.Cells.Item(105, 9).Value = 2210.3333  # BSM!I105 was 2224.2666
$ws.Cells.Item(105, 11).Value = 2210.3333  # BSM!K105 was 2224.2666
$ws.Cells.Item(105, 13).Value = -463.3332999999998  # BSM!M105 was -477.2665999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 4971.74  # BSM!H134 was 4814.423
$ws.Cells.Item(134, 9).Value = 2286.5588  # BSM!I134 was 2208.5
$ws.Cells.Item(134, 11).Value = 6859.676399999999  # BSM!K134 was 6625.5
$ws.Cells.Item(134, 13).Value = -4324.676399999999  # BSM!M134 was -4090.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5815.25  # CRP!H31 was 5891.7144
$ws.Cells.Item(31, 9).Value = 2431.279  # CRP!I31 was 2465.4048
$ws.Cells.Item(31, 11).Value = 2431.279  # CRP!K31 was 2465.4048
$ws.Cells.Item(31, 13).Value = -2136.279  # CRP!M31 was -2170.4048

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 5815.25  # CRP!H34 was 5891.7144
$ws.Cells.Item(34, 9).Value = 2431.279  # CRP!I34 was 2465.4048
$ws.Cells.Item(34, 11).Value = 2431.279  # CRP!K34 was 2465.4048
$ws.Cells.Item(34, 13).Value = -2229.279  # CRP!M34 was -2263.4048

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(56, 8).Value = 44875  # CRP!H56 was 45375
$ws.Cells.Item(56, 10).Value = 65000  # CRP!J56 was 66000
$ws.Cells.Item(56, 12).Value = 65000  # CRP!L56 was 66000
$ws.Cells.Item(56, 14).Value = -66690  # CRP!N56 was -67690

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(81, 8).Value = 41844  # CRP!H81 was 60000
$ws.Cells.Item(81, 10).Value = 41844  # CRP!J81 was 60000
$ws.Cells.Item(81, 12).Value = 41844  # CRP!L81 was 60000
$ws.Cells.Item(81, 14).Value = -43840  # CRP!N81 was -61996

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(84, 8).Value = 41844  # CRP!H84 was 60000
$ws.Cells.Item(84, 10).Value = 41844  # CRP!J84 was 60000
$ws.Cells.Item(84, 12).Value = 125532  # CRP!L84 was 180000
$ws.Cells.Item(84, 14).Value = -135516  # CRP!N84 was -189984

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 3646.8235  # CRP!H122 was 3699.75
$ws.Cells.Item(122, 9).Value = 3611.0833  # CRP!I122 was 3684.818
$ws.Cells.Item(122, 11).Value = 10833.2499  # CRP!K122 was 11054.454
$ws.Cells.Item(122, 13).Value = -8383.249899999999  # CRP!M122 was -8604.454000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 7692993.5  # CUL!H107 was 8000703
$ws.Cells.Item(107, 10).Value = 9524512  # CUL!J107 was 10000724
$ws.Cells.Item(107, 12).Value = 28573536  # CUL!L107 was 30002172
$ws.Cells.Item(107, 14).Value = -28577376  # CUL!N107 was -30006012

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(57, 8).Value = 57097.3  # GSM!H57 was 54633.637
$ws.Cells.Item(57, 10).Value = 61108.11  # GSM!J57 was 57997
$ws.Cells.Item(57, 12).Value = 61108.11  # GSM!L57 was 57997
$ws.Cells.Item(57, 14).Value = -62748.11  # GSM!N57 was -59637

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2461.6128  # GSM!H80 was 2415.9062
$ws.Cells.Item(80, 10).Value = 2355.4285  # GSM!J80 was 2185.875
$ws.Cells.Item(80, 12).Value = 2355.4285  # GSM!L80 was 2185.875
$ws.Cells.Item(80, 14).Value = -4351.4285  # GSM!N80 was -4181.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 2461.6128  # GSM!H83 was 2415.9062
$ws.Cells.Item(83, 10).Value = 2355.4285  # GSM!J83 was 2185.875
$ws.Cells.Item(83, 12).Value = 11777.1425  # GSM!L83 was 10929.375
$ws.Cells.Item(83, 14).Value = -21761.1425  # GSM!N83 was -20913.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 25814.889  # GSM!H122 was 28729.275
$ws.Cells.Item(122, 9).Value = 32338.295  # GSM!I122 was 36373.434
$ws.Cells.Item(122, 10).Value = 5651.636  # GSM!J122 was 5796.8
$ws.Cells.Item(122, 11).Value = 97014.88499999999  # GSM!K122 was 109120.302
$ws.Cells.Item(122, 12).Value = 16954.908  # GSM!L122 was 17390.4
$ws.Cells.Item(122, 13).Value = -94564.88499999999  # GSM!M122 was -106670.302
$ws.Cells.Item(122, 14).Value = -21854.908  # GSM!N122 was -22290.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 6977.647  # LTW!H22 was 6977.706
$ws.Cells.Item(22, 9).Value = 2649.6  # LTW!I22 was 1975.1428
$ws.Cells.Item(22, 10).Value = 8781  # LTW!J22 was 10479.5
$ws.Cells.Item(22, 11).Value = 2649.6  # LTW!K22 was 1975.1428
$ws.Cells.Item(22, 12).Value = 8781  # LTW!L22 was 10479.5
$ws.Cells.Item(22, 13).Value = -2354.6  # LTW!M22 was -1680.1428
$ws.Cells.Item(22, 14).Value = -9371  # LTW!N22 was -11069.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 6977.647  # LTW!H27 was 6977.706
$ws.Cells.Item(27, 9).Value = 2649.6  # LTW!I27 was 1975.1428
$ws.Cells.Item(27, 10).Value = 8781  # LTW!J27 was 10479.5
$ws.Cells.Item(27, 11).Value = 2649.6  # LTW!K27 was 1975.1428
$ws.Cells.Item(27, 12).Value = 8781  # LTW!L27 was 10479.5
$ws.Cells.Item(27, 13).Value = -2542.6  # LTW!M27 was -1868.1428
$ws.Cells.Item(27, 14).Value = -8995  # LTW!N27 was -10693.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1635.55  # LTW!H46 was 2001.4
$ws.Cells.Item(46, 9).Value = 889.7143  # LTW!I46 was 888.5
$ws.Cells.Item(46, 10).Value = 2037.1538  # LTW!J46 was 2525.1177
$ws.Cells.Item(46, 11).Value = 889.7143  # LTW!K46 was 888.5
$ws.Cells.Item(46, 12).Value = 2037.1538  # LTW!L46 was 2525.1177
$ws.Cells.Item(46, 13).Value = -701.7143  # LTW!M46 was -700.5
$ws.Cells.Item(46, 14).Value = -2413.1538  # LTW!N46 was -2901.1177

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 58823892  # LTW!H55 was 66667064
$ws.Cells.Item(55, 9).Value = 111111190  # LTW!I55 was 142857230
$ws.Cells.Item(55, 11).Value = 111111190  # LTW!K55 was 142857230
$ws.Cells.Item(55, 13).Value = -111111017  # LTW!M55 was -142857057

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 3007.5625  # LTW!H68 was 2789.111
$ws.Cells.Item(68, 9).Value = 1193  # LTW!I68 was 1169.6923
$ws.Cells.Item(68, 11).Value = 1193  # LTW!K68 was 1169.6923
$ws.Cells.Item(68, 13).Value = -444  # LTW!M68 was -420.6922999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value = 3007.5625  # LTW!H71 was 2789.111
$ws.Cells.Item(71, 9).Value = 1193  # LTW!I71 was 1169.6923
$ws.Cells.Item(71, 11).Value = 5965  # LTW!K71 was 5848.461499999999
$ws.Cells.Item(71, 13).Value = -2221  # LTW!M71 was -2104.461499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 3396.7222  # LTW!H82 was 3511.2354
$ws.Cells.Item(82, 9).Value = 2691.6667  # LTW!I82 was 2940
$ws.Cells.Item(82, 11).Value = 2691.6667  # LTW!K82 was 2940
$ws.Cells.Item(82, 13).Value = -2330.6667  # LTW!M82 was -2579

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(85, 8).Value = 3396.7222  # LTW!H85 was 3511.2354
$ws.Cells.Item(85, 9).Value = 2691.6667  # LTW!I85 was 2940
$ws.Cells.Item(85, 11).Value = 2691.6667  # LTW!K85 was 2940
$ws.Cells.Item(85, 13).Value = -1443.6667  # LTW!M85 was -1692

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 2792.7827  # LTW!H93 was 2874.3635
$ws.Cells.Item(93, 10).Value = 3333.3333  # LTW!J93 was 3800.4
$ws.Cells.Item(93, 12).Value = 3333.3333  # LTW!L93 was 3800.4
$ws.Cells.Item(93, 14).Value = -5829.3333  # LTW!N93 was -6296.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3672.121  # LTW!H122 was 3743.125
$ws.Cells.Item(122, 9).Value = 2617.5454  # LTW!I122 was 2675.524
$ws.Cells.Item(122, 11).Value = 7852.6362  # LTW!K122 was 8026.572
$ws.Cells.Item(122, 13).Value = -5402.6362  # LTW!M122 was -5576.572

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 6125.1665  # WVR!H62 was 5042.875
$ws.Cells.Item(62, 9).Value = 4238  # WVR!I62 was 3424
$ws.Cells.Item(62, 11).Value = 4238  # WVR!K62 was 3424
$ws.Cells.Item(62, 13).Value = -3614  # WVR!M62 was -2800

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(65, 8).Value = 6125.1665  # WVR!H65 was 5042.875
$ws.Cells.Item(65, 9).Value = 4238  # WVR!I65 was 3424
$ws.Cells.Item(65, 11).Value = 21190  # WVR!K65 was 17120
$ws.Cells.Item(65, 13).Value = -18070  # WVR!M65 was -14000

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 3542.121  # WVR!H126 was 3299.6943
$ws.Cells.Item(126, 9).Value = 2123.0435  # WVR!I126 was 1951.1154
$ws.Cells.Item(126, 11).Value = 6369.130500000001  # WVR!K126 was 5853.3462
$ws.Cells.Item(126, 13).Value = -3899.130500000001  # WVR!M126 was -3383.3462
